# Update the "dSF" column (F) values for several rows as part of a
# "repull data, push all data, mean calculation" refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    11 = 0
    14 = -3
    16 = -5
    17 = -4
    19 = -3
    24 = -1
    26 = -5
    29 = 1
    31 = -6
    32 = -8
    37 = -4
    38 = 0
    39 = 4
    43 = -5
    47 = -2
    50 = 4
    55 = -4
    56 = -9
    57 = -1
    61 = 3
    62 = -3
    64 = -1
    65 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
